# Rename the cash-flow category headers (C1:G1) to use underscores, and
# replace the old "Receipts_Ongoing Operations" header with "Repayments_loan".
# This adds a new cashflow line item (loan repayments) and normalizes the
# header naming convention used for the EUR conversion / 2H calculations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Receipts_Interest"
$ws.Range("D1").Value = "Receipts_Dividends"
$ws.Range("E1").Value = "New_Investments"
$ws.Range("F1").Value = "Development_Assets"
$ws.Range("G1").Value = "Repayments_loan"

# Resize the affected columns to fit their new (longer) header text, as
# Excel does automatically when column content changes width.
$ws.Range("C1:I10").EntireColumn.AutoFit()

# Leave the cursor where the user finished reviewing the update.
$ws.Range("F10").Select()
